$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "gulp" note to mention node js installation
$ws.Range("B60").Value = "галп + установка node js"

# New rows describing gulp npm packages (61-66)
$ws.Range("A61").Value = "https://www.npmjs.com/package/gulp-sass"
$ws.Range("A62").Value = "https://www.browsersync.io/"
$ws.Range("A63").Value = "https://www.npmjs.com/package/browser-sync"
$ws.Range("A64").Value = "https://www.npmjs.com/package/gulp-autoprefixer"

$ws.Range("B61").Value = "пакеты"
$ws.Range("B62").Value = "-"
$ws.Range("B63").Value = "-"
$ws.Range("B64").Value = "-"

$ws.Range("A65").Value = "https://www.npmjs.com/package/gulp-clean-css"
$ws.Range("A66").Value = "https://www.npmjs.com/package/gulp-rename"

$ws.Range("B65").Value = "-"
$ws.Range("B66").Value = "-"

# New section header "БЭМ" (row 67), styled/merged like the other section headers
$ws.Range("A59:B59").Copy()
$ws.Range("A67:B67").PasteSpecial(-4122)
$ws.Range("A67").Value = "БЭМ"
$ws.Range("A67:B67").Merge()

# Link to BEM methodology (row 68)
$ws.Range("A68").Value = "https://ru.bem.info/methodology/"

# Move selection to the newly added header row, matching the saved view
$ws.Range("A67:B67").Select()
